$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.425.32'
$ws.Range("E2").Value = '  +3.26%  '
$ws.Range("D3").Value = '3.065.42'
$ws.Range("E3").Value = '  +1.99%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''550.18'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").Value = '''140.15'
$ws.Range("E6").Value = '  +2.49%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.061.35'
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("D10").Value = '''6.52'
$ws.Range("E10").Value = '  +6.49%  '
$ws.Range("D11").Value = '''0.153'
$ws.Range("E11").Value = '  +1.95%  '
$ws.Range("D12").Value = '''0.456'
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("E13").Value = '  +2.34%  '
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").Value = '3.564.86'
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = '63.426.63'
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").Value = '3.063.59'
$ws.Range("E17").Value = '  +1.91%  '
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D20").Value = '''485.44'
$ws.Range("E20").Value = '  +3.34%  '
$ws.Range("D21").Value = '''13.79'
$ws.Range("E21").Value = '  +3.37%  '
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").Value = '''7.28'
$ws.Range("E23").Value = '  +3.86%  '
$ws.Range("D24").Value = '''81.01'
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("D25").Value = '''12.64'
$ws.Range("E25").Value = '  +4.36%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  +2.71%  '
$ws.Range("D28").Value = '''7.92'
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("E29").Value = '  +5.56%  '
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("E33").Value = '  +6.67%  '
$ws.Range("E34").Value = '  +3.15%  '
$ws.Range("D35").Value = '''55.63'
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  +0.96%  '
$ws.Range("D37").Value = '''466.36'
$ws.Range("E37").Value = '  +2.29%  '
$ws.Range("D38").Value = '''0.0823'
$ws.Range("E38").Value = '  +3.68%  '
$ws.Range("D39").Value = '''0.0398'
$ws.Range("E39").Value = '  +2.95%  '
$ws.Range("D40").Value = '3.055.87'
$ws.Range("E40").Value = '  -3.91%  '
$ws.Range("E41").Value = '  +0.97%  '
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("D43").Value = '''2.57'
$ws.Range("E43").Value = '  +2.15%  '
$ws.Range("D44").Value = '''28.20'
$ws.Range("E44").Value = '  +2.44%  '
$ws.Range("E45").Value = '  +3.57%  '
$ws.Range("E47").Value = '  +2.11%  '
$ws.Range("E48").Value = '  +1.61%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '''117.04'
$ws.Range("E49").Value = '  -2.56%  '
$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.0₃0512'
$ws.Range("E50").Value = '  +2.25%  '
$ws.Range("E51").Value = '  +2.50%  '
